# EngineUpgrade.xlsx update
# - Rebalance tier values (rows 2-5)
# - Add a new "ASL Thrust" row (row 7: label in A7, formula in B7)
# - Remove the stray helper formula that used to live in B8
# - Update the current sheet selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
$ws.Range("B2").Value = 1950
$ws.Range("C2").Value = 1250

# --- Row 3 ---------------------------------------------------------------
$ws.Range("B3").Value = 155000
$ws.Range("C3").Value = 1170

# --- Row 4 ---------------------------------------------------------------
$ws.Range("B4").Value = 295000
$ws.Range("C4").Value = 850

# --- Row 5 ---------------------------------------------------------------
$ws.Range("B5").Value = 12

# --- Row 7: new "ASL Thrust" entry ----------------------------------------
$ws.Range("A7").Value = "ASL Thrust"
$ws.Range("B7").Formula = "=C3/C2*B2"

# --- Row 8: drop the old helper formula -----------------------------------
$ws.Range("B8").ClearContents()

# --- Selection -------------------------------------------------------------
[void]$ws.Range("H7:L23").Select()
